$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a weekly "Poroto granado" price log in rows 2-57 (row 1 = headers).
# A new week's record is inserted at row 17, pushing the former rows 17-57 down
# to rows 18-58 (the dataset grows from A1:R57 to A1:R58).

# 1) Insert a blank row at position 18; this leaves row 17 untouched and shifts
#    every row from 18 downward (old 18..57) down by one, to 19..58.
$ws.Rows.Item(18).Insert()

# 2) The row that used to be 17 should now occupy row 18 (it shifted logically
#    one week later in the table), so copy the (still intact) row 17 values into
#    the newly-created blank row 18.
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(18, $c).Value2 = $ws.Cells.Item(17, $c).Value2
}

# 3) Overwrite row 17 with the new week's record: a new date and a new volume
#    (Volumen) figure; the remaining fields (market, region, category, prices,
#    unit, origin, $/Kg, Kg/Unidades, classification) match the row immediately
#    below it, so they are left as-is.
$ws.Cells.Item(17, 4).Value2 = 44575   # Fecha -> 2022-01-14
$ws.Cells.Item(17, 10).Value2 = 400    # Volumen
